{"js": "// Office.js (Word JavaScript API) script.\n// Applies the diff:\n//   1. Insert a new paragraph right after the title (Heading1) paragraph\n//      containing an empty leading run, a bold \"Meta description\" run and\n//      a normal run with the meta-description text.\n//   2. Near the end of the document, delete the paragraph that duplicated\n//      the page title (bold \"Play Atlantis Megaways Slot Game for Free |\n//      Review\") and replace the text of the following (italic) paragraph\n//      with the new AI image-generation prompt, keeping its italic run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1. Insert the \"Meta description\" paragraph after the H1 title -------\nconst titlePara = paragraphs.items[0];\nconst titleEnd = titlePara.getRange(\"End\");\n\nconst metaParaOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r/>' +\n  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +\n  '<w:r><w:t>: Read our Atlantis Megaways slot game review and play for free. Features, gameplay mechanics, betting range, and jackpot prizes explained.</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntitleEnd.insertOoxml(metaParaOoxml, \"After\");\nawait context.sync();\n\n// --- 2. Drop the duplicated bold title paragraph near the end, and swap --\n//        the italic paragraph's text for the new image prompt. -----------\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nconst duplicateTitlePara = paragraphs.items[count - 2];\nconst imagePromptPara = paragraphs.items[count - 1];\n\nduplicateTitlePara.delete();\nawait context.sync();\n\nconst oldDescriptionText =\n  \"Read our Atlantis Megaways slot game review and play for free. Features, gameplay mechanics, betting range, and jackpot prizes explained.\";\nconst matches = imagePromptPara.search(oldDescriptionText, { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nconst newImagePromptText =\n  \"Create a cartoon-style feature image for Atlantis Megaways that features a happy Maya warrior wearing glasses. The warrior should be positioned underwater among ruins of the lost city of Atlantis with sea creatures swimming around in the background. The image should incorporate the game's logo and feature vibrant colors that capture the adventurous and mysterious theme of the game. The image should also clearly convey the idea of winning cash prizes with a bubbly, celebratory vibe.\";\n\nmatches.items[0].insertText(newImagePromptText, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the diff:\n#   1. Insert a new paragraph right after the title (Heading1) paragraph\n#      containing an empty leading run, a bold \"Meta description\" run and\n#      a normal run with the meta-description text.\n#   2. Near the end of the document, delete the paragraph that duplicated\n#      the page title (bold \"Play Atlantis Megaways Slot Game for Free |\n#      Review\") and replace the text of the following (italic) paragraph\n#      with the new AI image-generation prompt, keeping its italic run.\n\n$d = $word.ActiveDocument\n\n# --- 1. Insert the \"Meta description\" paragraph after the H1 title -------\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$titleRange.Collapse(0)          # wdCollapseEnd\n$titleRange.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n$metaRange = $metaPara.Range\n$metaXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our Atlantis Megaways slot game review and play for free. Features, gameplay mechanics, betting range, and jackpot prizes explained.</w:t></w:r></w:p></w:body></w:document>'\n$metaRange.InsertXML($metaXml)\n\n# --- 2. Drop the duplicated bold title paragraph near the end, and swap --\n#        the italic paragraph's text for the new image prompt. -----------\n$count = $d.Paragraphs.Count\n$duplicateTitlePara = $d.Paragraphs.Item($count - 1)\n$duplicateTitlePara.Range.Delete()\n\n$count2 = $d.Paragraphs.Count\n$imagePromptPara = $d.Paragraphs.Item($count2)\n$imagePromptRange = $imagePromptPara.Range\n\n$oldDescriptionText = \"Read our Atlantis Megaways slot game review and play for free. Features, gameplay mechanics, betting range, and jackpot prizes explained.\"\n$newImagePromptText = \"Create a cartoon-style feature image for Atlantis Megaways that features a happy Maya warrior wearing glasses. The warrior should be positioned underwater among ruins of the lost city of Atlantis with sea creatures swimming around in the background. The image should incorporate the game's logo and feature vibrant colors that capture the adventurous and mysterious theme of the game. The image should also clearly convey the idea of winning cash prizes with a bubbly, celebratory vibe.\"\n\n# Locate the old description text and collapse the range onto it, then\n# assign .Text directly (NOT Find.Execute's Replacement, which runs the\n# text through Word's smart-quote AutoCorrect and would mangle the\n# straight apostrophe in \"game's\").\n$find = $imagePromptRange.Find\n$find.Text = $oldDescriptionText\n$find.Execute() | Out-Null\n$imagePromptRange.Text = $newImagePromptText\n"}
